# Applies "Completed parallel run and Extent Report implementation":
#  - Row 2: rename test case, reset invocation count, bump priority
#  - Row 4: rename test case (negative variant), reset Execute flag to No
#  - Row 5: rename description to "Testing Link Broken 2", reset invocation count
#  - Add new rows 6-9 for the searchFlights test cases
#  - Extend the Execute (Yes/No) list data validation down to row 9
#  - Move the active selection to D7:D9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : validateHotelIsSelected -> validateBrokenIsSelected ---
$ws.Range("A2").Value = "validateBrokenIsSelected"
$ws.Range("B2").Value = "Testing Link Broken"
$ws.Range("C2").Value = "No"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 5

# --- Row 3 : validateOfferesIsSelected (unchanged values, keep as-is) ---
$ws.Range("A3").Value = "validateOfferesIsSelected"
$ws.Range("B3").Value = "Testing Link Checkboxes"
$ws.Range("C3").Value = "No"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 2

# --- Row 4 : validateOfferesIsSelected2 -> validateOfferesIsSelectedNegative ---
$ws.Range("A4").Value = "validateOfferesIsSelectedNegative"
$ws.Range("B4").Value = "Testing Link Checkboxes Negative"
$ws.Range("C4").Value = "No"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 3

# --- Row 5 : validateHotelIsSelected2 description updated ---
$ws.Range("A5").Value = "validateHotelIsSelected2"
$ws.Range("B5").Value = "Testing Link Broken 2"
$ws.Range("C5").Value = "No"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4

# --- New rows 6-9 : searchFlights test cases ---
$ws.Range("A6").Value = "searchFlights"
$ws.Range("B6").Value = "Validating Search resrch result"
$ws.Range("C6").Value = "Yes"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

$ws.Range("A7").Value = "searchFlights2"
$ws.Range("B7").Value = "Validating Search resrch result"
$ws.Range("C7").Value = "No"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 1

$ws.Range("A8").Value = "searchFlights3"
$ws.Range("B8").Value = "Validating Search resrch result"
$ws.Range("C8").Value = "No"
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1

$ws.Range("A9").Value = "searchFlights4"
$ws.Range("B9").Value = "Validating Search resrch result"
$ws.Range("C9").Value = "No"
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# --- Extend the Yes/No list validation on Execute column down to row 9 ---
$ws.Range("C2:C9").Validation.Delete()
$ws.Range("C2:C9").Validation.Add(3, 1, 1, '"Yes,No"')
$ws.Range("C2:C9").Validation.InCellDropdown = $true
$ws.Range("C2:C9").Validation.IgnoreBlank = $true
$ws.Range("C2:C9").Validation.ShowInput = $true
$ws.Range("C2:C9").Validation.ShowError = $true

# --- Update the selection to match the saved view state ---
$ws.Range("D7:D9").Select()
